# Update "想去人数" (F column) figures across sheets, reflecting refreshed
# scrape numbers for the gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value  = 373
$wsExpo.Range("F8").Value  = 80
$wsExpo.Range("F9").Value  = 574
$wsExpo.Range("F12").Value = 1205
$wsExpo.Range("F15").Value = 1068
$wsExpo.Range("F16").Value = 421
$wsExpo.Range("F17").Value = 6769
$wsExpo.Range("F21").Value = 7695
$wsExpo.Range("F24").Value = 3790
$wsExpo.Range("F26").Value = 2196
$wsExpo.Range("F33").Value = 245
$wsExpo.Range("F36").Value = 1837
$wsExpo.Range("F40").Value = 517
$wsExpo.Range("F42").Value = 1296
$wsExpo.Range("F44").Value = 1958

# --- 本地生活 (Local life) sheet ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 1248

# --- 全部类型 (All types / combined) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 1248
$wsAll.Range("F7").Value  = 373
$wsAll.Range("F8").Value  = 80
$wsAll.Range("F9").Value  = 574
$wsAll.Range("F12").Value = 1205
$wsAll.Range("F15").Value = 1068
$wsAll.Range("F16").Value = 421
$wsAll.Range("F17").Value = 6769
$wsAll.Range("F21").Value = 7695
$wsAll.Range("F24").Value = 3791
$wsAll.Range("F26").Value = 2196
$wsAll.Range("F35").Value = 245
$wsAll.Range("F36").Value = 1837
$wsAll.Range("F40").Value = 517
$wsAll.Range("F43").Value = 1296
$wsAll.Range("F45").Value = 1958
